$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new match data
$ws.Range("A2").Value = "CvGRMovR"
$ws.Range("C2").Value = "16:30"
$ws.Range("D2").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E2").Value = "Chico"
$ws.Range("F2").Value = "Bucaramanga"

$ws.Range("G2").Value = 2.8
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 1.83
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 6
$ws.Range("AD2").Value = 6
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 201
$ws.Range("AH2").Value = 6.5
$ws.Range("AI2").Value = 12
$ws.Range("AJ2").Value = 12
$ws.Range("AK2").Value = 29
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 17
$ws.Range("AP2").Value = 34
$ws.Range("AQ2").Value = 51
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 351
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 4.5
$ws.Range("AX2").Value = 17
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 51
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 351
$ws.Range("BC2").Value = 126
$ws.Range("BD2").Value = 126
